$wb = $excel.ActiveWorkbook

# --- Reviews_AdvancedFilters: rework the Keyword/Tag columns ---
$wsAdv = $wb.Worksheets.Item("Reviews_AdvancedFilters")
$wsAdv.Range("E1").Value = "Tag"
$wsAdv.Range("F1:F2").Delete()
$wsAdv.Range("C2").Value = "1,2,3,Recommended,Not Recommended,No Rating"
$wsAdv.Range("E2").Value = "null"
$wsAdv.Columns("C:C").AutoFit()
$wsAdv.Range("D2").Select()

# --- Insert the new "Other_Filters" sheet before "Sentiment_Filters" ---
$wsSentiment = $wb.Worksheets.Item("Sentiment_Filters")
$wsOther = $wb.Worksheets.Add($wsSentiment)
$wsOther.Name = "Other_Filters"
$wsOther.Range("A1").Value = "Keywords"
$wsOther.Range("A2").Value = "and"
$wsOther.Range("A2").Select()
$wsOther.Activate()

Write-Output "done"
